$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rows 3-10 -- D, L, M columns become plain value 18 (drop formulas) ---
foreach ($r in 3..10) {
    $ws.Range("D$r").Value = 18
    $ws.Range("L$r").Value = 18
    $ws.Range("M$r").Value = 18
}

# --- Step 2: capture the threaded comment text anchored on the old C12 (kvæg_ringkanal row) ---
$oldComment = $ws.Range("C12").CommentThreaded
$commentText = $oldComment.Text()

# --- Step 3: delete row 11 (svin_gylle) -- everything below shifts up by one row ---
$ws.Rows("11").Delete()

# --- Step 4: rows 11-15 (post-shift) -- D and L columns become 29 ---
foreach ($r in 11..15) {
    $ws.Range("D$r").Value = 29
    $ws.Range("L$r").Value = 29
}

# --- Step 5: delete the trailing row (kvæg_gylle), now row 16 ---
$ws.Rows("16").Delete()

# --- Step 6: move the threaded comment from (stale) C12 to C11 ---
$oldComment.Delete()
$ws.Range("C11").AddCommentThreaded($commentText) | Out-Null

# --- Step 7: column widths (best effort -- engine quantizes to 1/6 px steps) ---
$ws.Columns("K").ColumnWidth = 14.0
$ws.Columns("L").ColumnWidth = 18.5

# --- Step 8: selection ---
$ws.Range("M3:M10").Select() | Out-Null
